$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the whole data block (Position 1..9 columns B..J, Move 1..5 rows 2..6)
# to their initial value of 0.
$ws.Range("B2:J6").Value = 0.0

# The sheet already has an existing styled cell (D2, style index 1) in this
# block; copy its formatting onto the whole range so every newly populated
# cell keeps a consistent style instead of picking up the workbook default.
$ws.Range("D2").Copy()
$ws.Range("B2:J6").PasteSpecial(-4122)
